# Remove the "Dependency_Type_Descrip" / "Business description..." column
# (old column F) from the diagram sheet, shifting the remaining columns
# (Dependency_Descrip, Organization) left by one. This also drops the now
# unused "Dependency_Type_Descrip" / "Business description..." shared
# strings automatically, since nothing references them anymore.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire column F (this shifts G->F and H->G, and updates the
# sheet dimension accordingly).
$ws.Columns("F").Delete()

# Re-create the selection state shown by the source file: the whole of
# (the new) column F selected, scrolled back to the top of the sheet.
$ws.Activate()
$ws.Range("A1").Select()
$ws.Columns("F").Select()
